# Weekly update: insert a new price record for "Vega Monumental Concepción - Zanahoria"
# right before the existing row 396, shifting the remaining historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 396 (existing rows 396:442 shift down to 397:443)
$ws.Rows.Item(396).Insert()

# Populate the newly inserted row with the new week's record
$ws.Range("A396").Value = 11
$ws.Range("B396").Value = "Vega Monumental Concepción"
$ws.Range("C396").Value = "Bíobío"
$ws.Range("D396").Value = 45142
$ws.Range("E396").Value = 8
$ws.Range("F396").Value = 100114013
$ws.Range("G396").Value = "Zanahoria"
$ws.Range("H396").Value = "Sin especificar"
$ws.Range("I396").Value = "Primera"
$ws.Range("J396").Value = 150
$ws.Range("K396").Value = 4000
$ws.Range("L396").Value = 4000
$ws.Range("M396").Value = 4000
$ws.Range("N396").Value = "`$/saco 20 kilos"
$ws.Range("O396").Value = "Región de La Araucanía"
$ws.Range("P396").Value = 200
$ws.Range("Q396").Value = 20
$ws.Range("R396").Value = "Hortaliza"
